$d = $word.ActiveDocument

# Locate the paragraph that holds the '{m:userdoc 'zone1'}' user-doc field.
# It is currently stored as two runs: "{m" and ":userdoc 'zone1'}".
# We need to split it into four runs: "{", "m", ":userdoc 'zone1'", "}"
# (mirrors the split already used for the "{m:" / "enduserdoc}" field
# later in the document), so the field rewriter can address each token
# (open brace, "m", the doc call, close brace) independently.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "{m:userdoc*zone1*}*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start
    $fieldText = "{m:userdoc 'zone1'}"

    # Offsets (relative to paragraph start) of the two new split points:
    #   after "{"                -> offset 1
    #   before the closing "}"   -> offset (len - 1)
    $afterBrace = $paraStart + 1
    $beforeCloseBrace = $paraStart + $fieldText.Length - 1

    # Splitting a run without changing any character formatting: adding
    # and immediately removing a bookmark at the boundary forces Word to
    # break the run in two, with no residual rPr on either half.

    # Split "{m" -> "{" | "m"
    $d.Bookmarks.Add("zzTmpSplit1", $d.Range($afterBrace, $afterBrace))
    $d.Bookmarks("zzTmpSplit1").Delete()

    # Split ":userdoc 'zone1'}" -> ":userdoc 'zone1'" | "}"
    $d.Bookmarks.Add("zzTmpSplit2", $d.Range($beforeCloseBrace, $beforeCloseBrace))
    $d.Bookmarks("zzTmpSplit2").Delete()
}
